$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 2).Value = 'Bitcoin'
$ws.Cells.Item(2, 3).Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
Set-TextValue $ws.Cells.Item(2, 4) '30.316.46'
Set-TextValue $ws.Cells.Item(2, 5) '  -0.14%  '

$ws.Cells.Item(3, 2).Value = 'Ethereum'
$ws.Cells.Item(3, 3).Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
Set-TextValue $ws.Cells.Item(3, 4) '1.872.12'
Set-TextValue $ws.Cells.Item(3, 5) '  +0.27%  '

$ws.Cells.Item(4, 2).Value = 'TetherUSD'
$ws.Cells.Item(4, 3).Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
Set-TextValue $ws.Cells.Item(4, 4) '0.9998'
Set-TextValue $ws.Cells.Item(4, 5) '  -0.08%  '

$ws.Cells.Item(5, 2).Value = 'BNB'
$ws.Cells.Item(5, 3).Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextValue $ws.Cells.Item(5, 4) '244.98'
Set-TextValue $ws.Cells.Item(5, 5) '  +4.45%  '

$ws.Cells.Item(6, 2).Value = 'USDC'
$ws.Cells.Item(6, 3).Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextValue $ws.Cells.Item(6, 4) '0.9999'
Set-TextValue $ws.Cells.Item(6, 5) '  -0.03%  '

$ws.Cells.Item(7, 2).Value = 'XRP'
$ws.Cells.Item(7, 3).Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextValue $ws.Cells.Item(7, 4) '0.4723'
Set-TextValue $ws.Cells.Item(7, 5) '  +0.27%  '

$ws.Cells.Item(8, 2).Value = 'Cardano'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue $ws.Cells.Item(8, 4) '0.2880'
Set-TextValue $ws.Cells.Item(8, 5) '  +0.48%  '

$ws.Cells.Item(9, 2).Value = 'Dogecoin'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue $ws.Cells.Item(9, 4) '0.06483'
Set-TextValue $ws.Cells.Item(9, 5) '  -1.23%  '

$ws.Cells.Item(10, 2).Value = 'Solana'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue $ws.Cells.Item(10, 4) '21.14'
Set-TextValue $ws.Cells.Item(10, 5) '  -1.15%  '

$ws.Cells.Item(11, 2).Value = 'TRON'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Cells.Item(11, 4) '0.07758'
Set-TextValue $ws.Cells.Item(11, 5) '  -1.52%  '

$ws.Cells.Item(12, 2).Value = 'WrappedEther'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Cells.Item(12, 4) '1.867.56'
Set-TextValue $ws.Cells.Item(12, 5) '  +0.03%  '

$ws.Cells.Item(13, 2).Value = 'Litecoin'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Cells.Item(13, 4) '95.25'
Set-TextValue $ws.Cells.Item(13, 5) '  -1.73%  '

$ws.Cells.Item(14, 2).Value = 'Polygon'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Cells.Item(14, 4) '0.7151'
Set-TextValue $ws.Cells.Item(14, 5) '  +3.26%  '

$ws.Cells.Item(15, 2).Value = 'Polkadot'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Cells.Item(15, 4) '5.110'
Set-TextValue $ws.Cells.Item(15, 5) '  +0.14%  '

$ws.Cells.Item(16, 2).Value = 'BitcoinCash'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Cells.Item(16, 4) '276.42'
Set-TextValue $ws.Cells.Item(16, 5) '  +3.03%  '

$ws.Cells.Item(17, 2).Value = 'WrappedBTC'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Cells.Item(17, 4) '30.308.66'
Set-TextValue $ws.Cells.Item(17, 5) '  +0.04%  '

$ws.Cells.Item(18, 2).Value = 'Avalanche'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Cells.Item(18, 4) '13.39'
Set-TextValue $ws.Cells.Item(18, 5) '  -3.80%  '

$ws.Cells.Item(19, 2).Value = 'ShibaInu'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Cells.Item(19, 4) '0.000007581'
Set-TextValue $ws.Cells.Item(19, 5) '  -0.87%  '

$ws.Cells.Item(20, 2).Value = 'Dai'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Cells.Item(20, 4) '1.001'
Set-TextValue $ws.Cells.Item(20, 5) '  +0.08%  '

$ws.Cells.Item(21, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Cells.Item(21, 4) '2.138.90'
Set-TextValue $ws.Cells.Item(21, 5) '  +1.13%  '

$ws.Cells.Item(22, 2).Value = 'BinanceUSD'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Cells.Item(22, 4) '0.9996'
Set-TextValue $ws.Cells.Item(22, 5) '  -0.10%  '

$ws.Cells.Item(23, 2).Value = 'Uniswap'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Cells.Item(23, 4) '5.231'
Set-TextValue $ws.Cells.Item(23, 5) '  +0.05%  '

$ws.Cells.Item(24, 2).Value = 'Chainlink'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Cells.Item(24, 4) '6.165'
Set-TextValue $ws.Cells.Item(24, 5) '  -0.22%  '

$ws.Cells.Item(25, 2).Value = 'Cosmos'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Cells.Item(25, 4) '9.283'
Set-TextValue $ws.Cells.Item(25, 5) '  -1.27%  '

$ws.Cells.Item(26, 2).Value = 'Monero'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Cells.Item(26, 4) '165.53'
Set-TextValue $ws.Cells.Item(26, 5) '  -0.96%  '

$ws.Cells.Item(27, 2).Value = 'EthereumClassic'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Cells.Item(27, 4) '18.92'
Set-TextValue $ws.Cells.Item(27, 5) '  +0.20%  '

$ws.Cells.Item(28, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Cells.Item(28, 4) '1.911'
Set-TextValue $ws.Cells.Item(28, 5) '  -1.78%  '

$ws.Cells.Item(29, 2).Value = 'Toncoin'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Cells.Item(29, 4) '1.378'
Set-TextValue $ws.Cells.Item(29, 5) '  +1.19%  '

$ws.Cells.Item(30, 2).Value = 'Stellar'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Cells.Item(30, 4) '0.09873'
Set-TextValue $ws.Cells.Item(30, 5) '  -0.43%  '

$ws.Cells.Item(31, 2).Value = 'PancakeSwap'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Cells.Item(31, 4) '1.519'
Set-TextValue $ws.Cells.Item(31, 5) '  +4.06%  '

$ws.Cells.Item(32, 2).Value = 'Filecoin'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Cells.Item(32, 4) '4.270'
Set-TextValue $ws.Cells.Item(32, 5) '  -2.36%  '

$ws.Cells.Item(33, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Cells.Item(33, 4) '4.038'
Set-TextValue $ws.Cells.Item(33, 5) '  -0.41%  '

$ws.Cells.Item(34, 2).Value = 'Hedera'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Cells.Item(34, 4) '0.04767'
Set-TextValue $ws.Cells.Item(34, 5) '  +0.36%  '

$ws.Cells.Item(35, 2).Value = 'ARBITRUM'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Cells.Item(35, 4) '1.123'
Set-TextValue $ws.Cells.Item(35, 5) '  -0.96%  '

$ws.Cells.Item(36, 2).Value = 'ImmutableX'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Cells.Item(36, 4) '0.6958'
Set-TextValue $ws.Cells.Item(36, 5) '  -0.96%  '

$ws.Cells.Item(37, 2).Value = 'HuobiToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Cells.Item(37, 4) '2.711'
Set-TextValue $ws.Cells.Item(37, 5) '  -0.24%  '

$ws.Cells.Item(38, 2).Value = 'VeChain'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Cells.Item(38, 4) '0.01852'
Set-TextValue $ws.Cells.Item(38, 5) '  -1.22%  '

$ws.Cells.Item(39, 2).Value = 'MXToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Cells.Item(39, 4) '2.743'
Set-TextValue $ws.Cells.Item(39, 5) '  -1.87%  '

$ws.Cells.Item(40, 2).Value = 'FraxShare'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Cells.Item(40, 4) '6.315'
Set-TextValue $ws.Cells.Item(40, 5) '  +0.11%  '

$ws.Cells.Item(41, 2).Value = 'Aave'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Cells.Item(41, 4) '70.76'
Set-TextValue $ws.Cells.Item(41, 5) '  -3.34%  '

$ws.Cells.Item(42, 2).Value = 'RenderToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Cells.Item(42, 4) '1.917'
Set-TextValue $ws.Cells.Item(42, 5) '  -1.76%  '

$ws.Cells.Item(43, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Cells.Item(43, 4) '0.8453'
Set-TextValue $ws.Cells.Item(43, 5) '  +0.27%  '

$ws.Cells.Item(44, 2).Value = 'PaxDollar'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws.Cells.Item(44, 4) '0.9997'
Set-TextValue $ws.Cells.Item(44, 5) '  -0.08%  '

$ws.Cells.Item(45, 2).Value = 'TheSandbox'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Cells.Item(45, 4) '0.4121'
Set-TextValue $ws.Cells.Item(45, 5) '  -1.41%  '

$ws.Cells.Item(46, 2).Value = 'Quant'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Cells.Item(46, 4) '101.89'
Set-TextValue $ws.Cells.Item(46, 5) '  -1.01%  '

$ws.Cells.Item(47, 2).Value = 'EnergySwap'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Cells.Item(47, 4) '9.294'
Set-TextValue $ws.Cells.Item(47, 5) '  +1.65%  '

$ws.Cells.Item(48, 2).Value = 'Aptos'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Cells.Item(48, 4) '7.115'
Set-TextValue $ws.Cells.Item(48, 5) '  -0.11%  '

$ws.Cells.Item(49, 2).Value = 'Elrond'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue $ws.Cells.Item(49, 4) '35.33'
Set-TextValue $ws.Cells.Item(49, 5) '  +2.39%  '

$ws.Cells.Item(50, 2).Value = 'Maker'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Cells.Item(50, 4) '921.55'
Set-TextValue $ws.Cells.Item(50, 5) '  -4.96%  '

$ws.Cells.Item(51, 2).Value = 'Cronos'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Cells.Item(51, 4) '0.05570'
Set-TextValue $ws.Cells.Item(51, 5) '  -1.90%  '
